# The underlying source re-sorted the sightings table. In terms of the
# worksheet this amounts to swapping the full content of row 9 with row 11,
# and the full content of row 10 with row 12 - including a handful of
# otherwise-blank "placeholder" cells (J/N/AC/AF) that only exist on row 11
# before the edit and need to move to row 9.
#
# Only cells that actually change value are touched (re-assigning a cell to
# the exact value it already holds is avoided) so that everything else in
# the sheet is left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by index) whose value differs between row 9 and row 11, and
# between row 10 and row 12 - i.e. the payload columns of one sighting
# record. K (11) is intentionally left out of the 9/11 swap because both
# rows already hold the same "teleomorf" value there.
$swapCols_9_11  = @(1, 2, 4, 5, 6, 7, 8, 17, 18)            # A B D E F G H Q R
$swapCols_10_12 = @(1, 2, 4, 5, 6, 7, 8, 11, 17, 18)        # A B D E F G H K Q R

function Get-RowValues($sheet, $row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $sheet.Cells.Item($row, $c).Value()
    }
    return $vals
}

function Set-RowValues($sheet, $row, $vals, $cols) {
    foreach ($c in $cols) {
        $sheet.Cells.Item($row, $c).Value = $vals[$c]
    }
}

$row9  = Get-RowValues $ws 9  $swapCols_9_11
$row11 = Get-RowValues $ws 11 $swapCols_9_11
Set-RowValues $ws 9  $row11 $swapCols_9_11
Set-RowValues $ws 11 $row9  $swapCols_9_11

$row10 = Get-RowValues $ws 10 $swapCols_10_12
$row12 = Get-RowValues $ws 12 $swapCols_10_12
Set-RowValues $ws 10 $row12 $swapCols_10_12
Set-RowValues $ws 12 $row10 $swapCols_10_12

# Row 9 gains the public-comment text that used to live on row 11.
$ws.Cells.Item(9, 29).Value = "Mörkröda droppar på hattöversidan. Smak besk efter ett långt tag, ej brännande. Köttet färgas mörkviolett med KOH."

# Row 9 also gains the three blank placeholder cells (J, N, AF) that used to
# exist only on row 11. Re-applying a cell's own (unchanged) number format
# is enough to make Excel materialise it as a real, present-but-empty cell.
foreach ($col in @(10, 14, 32)) {  # J, N, AF
    $cell = $ws.Cells.Item(9, $col)
    $cell.NumberFormat = $cell.NumberFormat
}

# Row 11 loses its public comment and its three blank placeholder cells -
# ClearContents() removes a cell's content; for already-empty placeholder
# cells that makes them disappear entirely, matching row 9's old shape.
$ws.Cells.Item(11, 29).ClearContents()  # AC11
$ws.Cells.Item(11, 10).ClearContents()  # J11
$ws.Cells.Item(11, 14).ClearContents()  # N11
$ws.Cells.Item(11, 32).ClearContents()  # AF11
